# Monte Carlo Type AI - data tweak: HP CUR / HP MAX values for rows 2 and 3
# plus refreshed selection / window chrome to match the latest save state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update HP CUR (H) / HP MAX (I) for Lugia (row 2) and Omastar (row 3) ---
$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 120
$ws.Range("H3").Value = 78
$ws.Range("I3").Value = 78

# --- Move the sheet's selection / active cell to I2 (also clears the old
#     scrolled-to-U1 top-left-cell state since the view resets on select) ---
$ws.Activate()
$ws.Range("I2").Select()

# --- Refresh the application / window geometry recorded for the workbook ---
try {
    $excel.Left = 240
    $excel.Top = 240
    $excel.Width = 25360
    $excel.Height = 15820
} catch {}

try {
    $win = $excel.ActiveWindow
    $win.Left = 240
    $win.Top = 240
    $win.Width = 25360
    $win.Height = 15820
} catch {}
